$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$rng = $ws.Range("D1:G$lastRow")

# Force the whole range to Text format first so numeric-looking codes like
# "110" / "111" round-trip as text (matching the source data) instead of
# being auto-converted to numbers by Excel's type inference.
$rng.NumberFormat = "@"

$vals = $rng.Value()

for ($i = 1; $i -le $lastRow; $i++) {
    $d = $vals[$i, 1]
    $e = $vals[$i, 2]
    $f = $vals[$i, 3]
    $g = $vals[$i, 4]

    # Swap D <-> E (category-name <-> group-name)
    $vals[$i, 1] = $e
    $vals[$i, 2] = $d

    # Swap F <-> G (group-code <-> category-code)
    $vals[$i, 3] = $g
    $vals[$i, 4] = $f
}

$rng.Value = $vals

# Restore the default "Normal" style so the cells don't carry a lingering
# custom number format now that the text values are safely in place.
$rng.Style = "Normal"
